$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# --- Crime Complaints table updates ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 22.222222222222
$ws.Range("L15").Value = 22.222222222222
$ws.Range("M15").Value = 37.5
$ws.Range("N15").Value = -21.428571428571
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 400
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 56
$ws.Range("K16").Value = 19.642857142857
$ws.Range("L16").Value = 4.6875
$ws.Range("M16").Value = -4.285714285714
$ws.Range("N16").Value = -85.434782608695
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 21.052631578947
$ws.Range("I17").Value = 129
$ws.Range("J17").Value = 92
$ws.Range("K17").Value = 40.217391304347
$ws.Range("L17").Value = 2.380952380952
$ws.Range("M17").Value = 101.5625
$ws.Range("N17").Value = -37.681159420289
$ws.Range("C18").Value = 4
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 260
$ws.Range("I18").Value = 75
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = 41.509433962264
$ws.Range("L18").Value = -21.875
$ws.Range("M18").Value = -33.62831858407
$ws.Range("N18").Value = -83.14606741573
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -38.775510204081
$ws.Range("I19").Value = 198
$ws.Range("J19").Value = 235
$ws.Range("K19").Value = -15.744680851063
$ws.Range("L19").Value = -45.753424657534
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 15.78947368421
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 11
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 66
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = 6.451612903225
$ws.Range("L20").Value = 17.857142857142
$ws.Range("M20").Value = 29.411764705882
$ws.Range("N20").Value = -84.360189573459
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 22.222222222222
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = 5.154639175257
$ws.Range("I21").Value = 546
$ws.Range("J21").Value = 508
$ws.Range("K21").Value = 7.480314960629
$ws.Range("L21").Value = -23.743016759776
$ws.Range("M21").Value = 23.529411764705
$ws.Range("N21").Value = -68.384481760277
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -20.833333333333
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -20.481927710843
$ws.Range("I24").Value = 378
$ws.Range("J24").Value = 463
$ws.Range("K24").Value = -18.358531317494
$ws.Range("L24").Value = -27.862595419847
$ws.Range("M24").Value = 24.342105263157
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 179
$ws.Range("J25").Value = 258
$ws.Range("K25").Value = -30.620155038759
$ws.Range("L25").Value = -42.811501597444
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -30
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 45.714285714285
$ws.Range("I26").Value = 200
$ws.Range("J26").Value = 182
$ws.Range("K26").Value = 9.890109890109
$ws.Range("L26").Value = 22.699386503067
$ws.Range("M26").Value = -19.354838709677
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = 7.142857142857
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 6
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E28").Value = -83.333333333333
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = -73.333333333333
$ws.Range("I28").Value = 23
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -36.111111111111
$ws.Range("L28").Value = -8
$ws.Range("M29").Value = -66.666666666666
$ws.Range("M30").Value = -60
$ws.Range("F33").Value = 2
$ws.Range("I33").Value = 3
$ws.Range("L33").Value = 200
